# Apply updated crypto price/volume data to worksheet cells.
# Cells in column D that look numeric must be forced to stay as text
# (matching the original inline-string / text cell type) by setting
# NumberFormat to "@" (Text) before assigning the value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.925.21"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "2.293.46"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.54"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.64"
$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").Value = "2.313.62"
$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("E10").Value = "  +2.62%  "

$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.11"
$ws.Range("E12").Value = "  +7.81%  "

$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.82"
$ws.Range("E14").Value = "  +3.96%  "

$ws.Range("D15").Value = "2.701.49"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").Value = "54.916.19"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("E17").Value = "  +1.66%  "

$ws.Range("D18").Value = "2.292.10"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.52"
$ws.Range("E19").Value = "  +2.43%  "

$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.39"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  +4.02%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.28"
$ws.Range("E24").Value = "  -2.67%  "

$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.151"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.97"
$ws.Range("E28").Value = "  +0.91%  "

$ws.Range("E29").Value = "  +3.19%  "

$ws.Range("D30").Value = "0.0₃0710"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("E32").Value = "  +4.84%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.06"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("E37").Value = "  +2.73%  "

$ws.Range("E38").Value = "  +3.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.81"
$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("E40").Value = "  +2.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.376"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "133.89"
$ws.Range("E42").Value = "  +6.11%  "

$ws.Range("E43").Value = "  +1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.93"
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "261.32"
$ws.Range("E45").Value = "  +7.92%  "

$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0913"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("E50").Value = "  +2.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.48"
$ws.Range("E51").Value = "  +0.60%  "
